# Updated Master Data excels
# - Remove the extra "Sheet1" helper/staging sheet (and its Table1)
# - Append the newly-defined UIN card/deactivation/reactivation and
#   registration-acknowledgement template master-data rows to
#   "master-template_type"

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Drop the old "Sheet1" staging/filter sheet (it only ever held a
#    filtered copy of the same table + Table1 autofilter definition).
# ---------------------------------------------------------------------
foreach ($sh in @($wb.Worksheets)) {
    if ($sh.Name -eq "Sheet1") {
        $sh.Delete()
    }
}

$ws = $wb.Worksheets.Item("master-template_type")

# ---------------------------------------------------------------------
# 2. Append the new template rows (code, descr, lang_code) - is_active,
#    cr_by and cr_dtimes are the same for every row in this sheet.
# ---------------------------------------------------------------------
$newRows = @(
    @('RPR_UIN_CARD_TEMPLATE','UIN card template','eng'),
    @('RPR_UIN_CARD_TEMPLATE','قالب بطاقة UIN','ara'),
    @('RPR_UIN_CARD_TEMPLATE','Modèle de carte UIN','fra'),
    @('RPR_UIN_DEAC_SMS','Template for UIN Deactivation SMS','eng'),
    @('RPR_UIN_DEAC_SMS','قالب لتعطيل UIN SMS','ara'),
    @('RPR_UIN_DEAC_SMS','Modèle pour SMS de désactivation UIN','fra'),
    @('RPR_UIN_DEAC_EMAIL','Template for UIN Deactivation Email','eng'),
    @('RPR_UIN_DEAC_EMAIL','قالب لإلغاء تنشيط البريد','ara'),
    @('RPR_UIN_DEAC_EMAIL','Modèle pour Email de désactivation UIN','fra'),
    @('RPR_UIN_REAC_SMS','Template for UIN Reactivate SMS','eng'),
    @('RPR_UIN_REAC_SMS','قالب لـ UIN تنشيط SMS','ara'),
    @('RPR_UIN_REAC_SMS','Modèle pour UIN Réactiver SMS','fra'),
    @('RPR_UIN_REAC_EMAIL','Template for UIN Reactivate Email','eng'),
    @('RPR_UIN_REAC_EMAIL','قالب لـ UIN تنشيط البريد','ara'),
    @('RPR_UIN_REAC_EMAIL','Modèle pour UIN Réactiver Email','fra'),
    @('reg-sms-notification','Registration Acknowledgement Template','eng'),
    @('reg-sms-notification','نموذج شكر التسجيل','ara'),
    @('reg-sms-notification','accusé de réception','fra'),
    @('reg-email-notification','Registration Acknowledgement Template','eng'),
    @('reg-email-notification','نموذج شكر التسجيل','ara'),
    @('reg-email-notification','accusé de réception','fra'),
    @('reg-ack-template-part1','Registration Acknowledgement Template - Part 1','eng'),
    @('reg-ack-template-part2','نموذج شكر التسجيل','ara'),
    @('reg-ack-template-part3','accusé de réception','fra'),
    @('reg-ack-template-part2','Registration Acknowledgement Template - Part 2','eng'),
    @('reg-ack-template-part3','نموذج شكر التسجيل','ara'),
    @('reg-ack-template-part4','accusé de réception','fra'),
    @('reg-ack-template-part3','Registration Acknowledgement Template - Part 3','eng'),
    @('reg-ack-template-part4','نموذج شكر التسجيل','ara'),
    @('reg-ack-template-part5','accusé de réception','fra')
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Leave the selection on the row below the new data, matching the
#    "select to end of sheet" pattern already used on this tab.
# ---------------------------------------------------------------------
$lastRow = $r
$ws.Range("A" + $lastRow + ":XFD1048576").Select()
